$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.415.68'
$ws.Range('E2').Value = '  +2.46%  '
$ws.Range('D3').Value = '3.379.84'
$ws.Range('E3').Value = '  +1.77%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''588.05'
$ws.Range('E5').Value = '  +1.64%  '
$ws.Range('D6').Value = '''179.48'
$ws.Range('E6').Value = '  +2.62%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('E8').Value = '  +1.38%  '
$ws.Range('E9').Value = '  +5.04%  '
$ws.Range('D10').Value = '''0.588'
$ws.Range('E10').Value = '  +2.25%  '
$ws.Range('D11').Value = '''48.32'
$ws.Range('E11').Value = '  +6.58%  '
$ws.Range('E12').Value = '  +3.18%  '
$ws.Range('D13').Value = '''700.28'
$ws.Range('E13').Value = '  +6.01%  '
$ws.Range('D14').Value = '3.938.62'
$ws.Range('E14').Value = '  +1.95%  '
$ws.Range('E15').Value = '  +1.81%  '
$ws.Range('D16').Value = '69.358.91'
$ws.Range('E16').Value = '  +2.55%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.384.35'
$ws.Range('E17').Value = '  +1.72%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = '''0.120'
$ws.Range('E19').Value = '  +1.42%  '
$ws.Range('D20').Value = '''11.31'
$ws.Range('E20').Value = '  +3.45%  '
$ws.Range('D21').Value = '''0.902'
$ws.Range('E21').Value = '  +1.97%  '
$ws.Range('D22').Value = '''5.51'
$ws.Range('E22').Value = '  +3.03%  '
$ws.Range('D23').Value = '''17.20'
$ws.Range('E23').Value = '  +1.80%  '
$ws.Range('D24').Value = '''101.63'
$ws.Range('E24').Value = '  +3.51%  '
$ws.Range('D25').Value = '''3.94'
$ws.Range('E25').Value = '  +2.69%  '
$ws.Range('E26').Value = '  +2.32%  '
$ws.Range('D27').Value = '''9.61'
$ws.Range('E27').Value = '  +4.05%  '
$ws.Range('D28').Value = '''33.49'
$ws.Range('E28').Value = '  +0.87%  '
$ws.Range('E29').Value = '  +2.84%  '
$ws.Range('E30').Value = '  -2.29%  '
$ws.Range('E31').Value = '  +2.18%  '
$ws.Range('D32').Value = '''553.94'
$ws.Range('E32').Value = '  -2.33%  '
$ws.Range('E33').Value = '  +1.72%  '
$ws.Range('B34').Value = 'dogwifhat'
$ws.Range('C34').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D34').Value = '''3.52'
$ws.Range('E34').Value = '  +8.47%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = '''58.26'
$ws.Range('E35').Value = '  +3.53%  '
$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').Value = '3.725.96'
$ws.Range('E36').Value = '  +1.62%  '
$ws.Range('D37').Value = '''0.999'
$ws.Range('E37').Value = '  -0.24%  '
$ws.Range('D38').Value = '''0.144'
$ws.Range('E38').Value = '  +10.64%  '
$ws.Range('D39').Value = '''34.95'
$ws.Range('E39').Value = '  +1.46%  '
$ws.Range('E40').Value = '  +3.61%  '
$ws.Range('E41').Value = '  +1.31%  '
$ws.Range('D42').Value = '0.0₃0684'
$ws.Range('E42').Value = '  +3.67%  '
$ws.Range('E43').Value = '  +2.11%  '
$ws.Range('E44').Value = '  +3.73%  '
$ws.Range('E45').Value = '  -2.41%  '
$ws.Range('E46').Value = '  +2.86%  '
$ws.Range('E47').Value = '  +1.76%  '
$ws.Range('E48').Value = '  -0.11%  '
$ws.Range('E49').Value = '  -1.20%  '
$ws.Range('D50').Value = '''132.13'
$ws.Range('E50').Value = '  +3.45%  '
$ws.Range('D51').Value = '''2.65'
$ws.Range('E51').Value = '  -0.72%  '
